# tut06/output/2001MM05.xlsx - "Updated solution for Tutorial 6"
#
# For each attendance date row:
#   - reformat the date in column A from dd/mm/yyyy to dd-mm-yyyy (stays plain text)
#   - refresh the Real/Duplicate/Invalid/Absent tally in columns D-H
#
# Excel auto-parses an unquoted "dd-mm-yyyy" string as a real date whenever the
# first segment is <= 12 (could be a month), which would turn the cell into a
# date serial instead of leaving it as text. To stop that, those cells are
# temporarily switched to a Text number format before the write, then the
# style is reset back to Normal afterwards so no stray formatting is left
# behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateText($cell, $text) {
    $rng = $ws.Range($cell)
    $day = [int]($text.Split("-")[0])
    if ($day -le 12) {
        $rng.NumberFormat = "@"
        $rng.Value = $text
        $rng.Style = "Normal"
    } else {
        $rng.Value = $text
    }
}

# row -> [date, D, E, F, G, H]
$rows = @(
    @(3,  "28-07-2022", 1, 0, 0, 1, 1),
    @(4,  "01-08-2022", 1, 1, 0, 0, 0),
    @(5,  "04-08-2022", 1, 1, 0, 0, 0),
    @(6,  "08-08-2022", 1, 1, 0, 0, 0),
    @(7,  "11-08-2022", 0, 0, 0, 0, 1),
    @(8,  "15-08-2022", 0, 0, 0, 0, 1),
    @(9,  "18-08-2022", 0, 0, 0, 0, 1),
    @(10, "22-08-2022", 0, 0, 0, 0, 1),
    @(11, "25-08-2022", 0, 0, 0, 0, 1),
    @(12, "29-08-2022", 0, 0, 0, 0, 1),
    @(13, "01-09-2022", 1, 1, 0, 0, 0),
    @(14, "05-09-2022", 0, 0, 0, 0, 1),
    @(15, "08-09-2022", 0, 0, 0, 0, 1),
    @(16, "12-09-2022", 0, 0, 0, 0, 1),
    @(17, "15-09-2022", 0, 0, 0, 0, 1),
    @(18, "19-09-2022", 0, 0, 0, 0, 1),
    @(19, "22-09-2022", 0, 0, 0, 0, 1),
    @(20, "26-09-2022", 0, 0, 0, 0, 1),
    @(21, "29-09-2022", 0, 0, 0, 0, 1)
)

# before-edit D:H baseline was uniformly 0,0,0,0,1 for every row - only write
# a cell back out when its target value actually differs from that baseline.
$baseline = @(0, 0, 0, 0, 1)
$cols = @("D", "E", "F", "G", "H")

foreach ($row in $rows) {
    $r = $row[0]
    $date = $row[1]

    Set-DateText "A$r" $date

    for ($i = 0; $i -lt 5; $i++) {
        $target = $row[2 + $i]
        if ($target -ne $baseline[$i]) {
            $ws.Range("$($cols[$i])$r").Value = $target
        }
    }
}
